$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the authoritative diff.
# Numeric-looking Price (D) values need NumberFormat "@" forced first so Excel
# keeps them as text (matching the source data which stores prices as text),
# instead of silently re-interpreting them as numbers.

$ws.Range("D2").Value = "68.286.50"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.598.25"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.03"
$ws.Range("E5").Value = "  +8.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "569.24"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.682"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "63.89"
$ws.Range("E10").Value = "  +14.60%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.21"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "4.168.61"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "3.598.49"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.25"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "68.139.92"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.24"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "405.57"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.39"
$ws.Range("E23").Value = "  +8.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.03"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.90"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.87"
$ws.Range("E27").Value = "  +6.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "734.39"
$ws.Range("E29").Value = "  +15.55%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.68"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.59"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.17"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.64"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.99"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  +5.81%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  +23.91%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "3.195.38"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.132"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("E45").Value = "  +10.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.10"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.56"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  -1.43%  "
